$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 11 new rows starting at row 153. This pushes the existing
#    "Development type" (153-161) and "Voluntary agreement" (162-163) blocks
#    down to 164-172 and 173-174 respectively, and Excel automatically
#    updates the existing merged cell ranges to match.
$ws.Rows("153:163").Insert()

# 2. Widen column D from 27 to 29 characters.
$ws.Columns.Item(4).ColumnWidth = 28.17

# 3. Populate the newly inserted rows (153-163) with the new
#    "Oil and gas permission types" module content.

# Row 153 - module header row
$ws.Range("A153").Value = "Oil and gas permission types"
$ws.Range("B153").Value = "Module for details about types of onshore oil and gas extraction permissions already received and applying for`n"
$ws.Range("C153").Value = "Oil and gas permission types[]"
$ws.Range("D153").Value = ""
$ws.Range("E153").Value = ""
$ws.Range("F153").Value = ""
$ws.Range("G153").Value = "List of permission types being applied for"
$ws.Range("H153").Value = "enum"
$ws.Range("I153").Value = "MUST"

# Row 154
$ws.Range("C154").Value = "Related permissions[]"
$ws.Range("D154").Value = "Reference"
$ws.Range("E154").Value = ""
$ws.Range("F154").Value = ""
$ws.Range("G154").Value = "The reference for the related application that permission was received for"
$ws.Range("H154").Value = "string"
$ws.Range("I154").Value = "MUST"

# Row 155
$ws.Range("C155").Value = "Related permissions[]"
$ws.Range("D155").Value = "Oil and gas permission type"
$ws.Range("E155").Value = ""
$ws.Range("F155").Value = ""
$ws.Range("G155").Value = "An oil and gas related permission type"
$ws.Range("H155").Value = "enum"
$ws.Range("I155").Value = "MUST"

# Row 156
$ws.Range("C156").Value = "Related permissions[]"
$ws.Range("D156").Value = "Decision date"
$ws.Range("E156").Value = ""
$ws.Range("F156").Value = ""
$ws.Range("G156").Value = "The date when the decision was made, in YYYY-MM-DD format"
$ws.Range("H156").Value = "string"
$ws.Range("I156").Value = "MUST"

# Row 157
$ws.Range("C157").Value = "Related permissions[]"
$ws.Range("D157").Value = "Condition number"
$ws.Range("E157").Value = ""
$ws.Range("F157").Value = ""
$ws.Range("G157").Value = "Number of any condition being breached"
$ws.Range("H157").Value = "string"
$ws.Range("I157").Value = "MAY"

# Row 158
$ws.Range("C158").Value = "Other details"
$ws.Range("D158").Value = ""
$ws.Range("E158").Value = ""
$ws.Range("F158").Value = ""
$ws.Range("G158").Value = "Explanation if other ground is selected"
$ws.Range("H158").Value = "string"
$ws.Range("I158").Value = "MAY"

# Row 159
$ws.Range("C159").Value = "Will consolidate permissions"
$ws.Range("D159").Value = ""
$ws.Range("E159").Value = ""
$ws.Range("F159").Value = ""
$ws.Range("G159").Value = "Is the applicant looking to consolidate permissions?"
$ws.Range("H159").Value = "boolean"
$ws.Range("I159").Value = "MUST"

# Row 160
$ws.Range("C160").Value = "Details"
$ws.Range("D160").Value = ""
$ws.Range("E160").Value = ""
$ws.Range("F160").Value = ""
$ws.Range("G160").Value = "Details about the consolidation or update of permissions"
$ws.Range("H160").Value = "string"
$ws.Range("I160").Value = "MAY"

# Row 161
$ws.Range("C161").Value = "Related proposals[]"
$ws.Range("D161").Value = "Reference"
$ws.Range("E161").Value = ""
$ws.Range("F161").Value = ""
$ws.Range("G161").Value = "The reference for the related application"
$ws.Range("H161").Value = "string"
$ws.Range("I161").Value = "MUST"

# Row 162
$ws.Range("C162").Value = "Related proposals[]"
$ws.Range("D162").Value = "Application type"
$ws.Range("E162").Value = ""
$ws.Range("F162").Value = ""
$ws.Range("G162").Value = "The type of planning application"
$ws.Range("H162").Value = "enum"
$ws.Range("I162").Value = "MUST"

# Row 163
$ws.Range("C163").Value = "Related proposals[]"
$ws.Range("D163").Value = "Decision date"
$ws.Range("E163").Value = ""
$ws.Range("F163").Value = ""
$ws.Range("G163").Value = "The date when the decision was made, in YYYY-MM-DD format"
$ws.Range("H163").Value = "string"
$ws.Range("I163").Value = "MUST"

# 4. Re-create the merges for the new module block (A153:A163, B153:B163).
$ws.Range("A153:A163").Merge()
$ws.Range("B153:B163").Merge()

Write-Output "done"
